# Swap the data of row 10 and row 11 for columns A, B, E, F, G, H, Q, R
# (the other columns are identical between the two rows, so no changes
# are needed there).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell10 = $ws.Range("${col}10")
    $cell11 = $ws.Range("${col}11")

    $val10 = $cell10.Value2
    $val11 = $cell11.Value2

    $cell10.Value = $val11
    $cell11.Value = $val10
}
